$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2857.1428
$ws.Range("J40").Value = 3500
$ws.Range("L40").Value = 3500
$ws.Range("N40").Value = -3850
$ws.Range("H70").Value = 110965.86
$ws.Range("I70").Value = 1163
$ws.Range("K70").Value = 3489
$ws.Range("M70").Value = -3219
$ws.Range("H73").Value = 110965.86
$ws.Range("I73").Value = 1163
$ws.Range("K73").Value = 3489
$ws.Range("M73").Value = -2553
$ws.Range("H132").Value = 1056.8422
$ws.Range("I132").Value = 1011.0909
$ws.Range("K132").Value = 3033.2727
$ws.Range("M132").Value = -503.2727
$ws.Range("H138").Value = 4384.0864
$ws.Range("J138").Value = 4776.16
$ws.Range("L138").Value = 14328.48
$ws.Range("N138").Value = -24608.48
$ws.Range("H141").Value = 5407.778
$ws.Range("I141").Value = 5310
$ws.Range("J141").Value = 5750
$ws.Range("K141").Value = 15930
$ws.Range("L141").Value = 17250
$ws.Range("M141").Value = -10750
$ws.Range("N141").Value = -27610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 751501.5
$ws.Range("I8").Value = 3000000
$ws.Range("K8").Value = 3000000
$ws.Range("M8").Value = -2999856
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents() | Out-Null
$ws.Range("H12").Value = 426
$ws.Range("I12").Value = 752
$ws.Range("K12").Value = 752
$ws.Range("M12").Value = -579
$ws.Range("H14").Value = 1603.25
$ws.Range("I14").Value = 1835.3334
$ws.Range("K14").Value = 1835.3334
$ws.Range("M14").Value = -1660.3334
$ws.Range("H16").Value = 5341.6665
$ws.Range("I16").Value = 512.5
$ws.Range("K16").Value = 512.5
$ws.Range("M16").Value = -225.5
$ws.Range("H21").Value = 657.8
$ws.Range("J21").Value = 758
$ws.Range("L21").Value = 758
$ws.Range("N21").Value = -1506
$ws.Range("H30").Value = 470
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 470
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 470
$ws.Range("M30").ClearContents() | Out-Null
$ws.Range("N30").Value = -770
$ws.Range("H32").Value = 16048.293
$ws.Range("I32").Value = 7276.625
$ws.Range("J32").Value = 26844.191
$ws.Range("K32").Value = 7276.625
$ws.Range("L32").Value = 26844.191
$ws.Range("M32").Value = -6989.625
$ws.Range("N32").Value = -27418.191
$ws.Range("H33").Value = 1800
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents() | Out-Null
$ws.Range("H132").Value = 1435.2894
$ws.Range("I132").Value = 1459.5555
$ws.Range("J132").Value = 998.5
$ws.Range("K132").Value = 4378.666499999999
$ws.Range("L132").Value = 2995.5
$ws.Range("M132").Value = -1848.666499999999
$ws.Range("N132").Value = -8055.5
$ws.Range("H138").Value = 85000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 85000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 85000
$ws.Range("M138").ClearContents() | Out-Null
$ws.Range("N138").Value = -95280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3184.675
$ws.Range("I105").Value = 2465.138
$ws.Range("J105").Value = 5081.636
$ws.Range("K105").Value = 2465.138
$ws.Range("L105").Value = 5081.636
$ws.Range("M105").Value = -718.1379999999999
$ws.Range("N105").Value = -8575.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 138.91667
$ws.Range("I7").Value = 227.85715
$ws.Range("J7").Value = 14.4
$ws.Range("K7").Value = 227.85715
$ws.Range("L7").Value = 14.4
$ws.Range("M7").Value = -114.85715
$ws.Range("N7").Value = -240.4
$ws.Range("H22").Value = 276.91666
$ws.Range("I22").Value = 279
$ws.Range("J22").Value = 275.42856
$ws.Range("K22").Value = 279
$ws.Range("L22").Value = 275.42856
$ws.Range("M22").Value = 71
$ws.Range("N22").Value = -975.4285600000001
$ws.Range("H31").Value = 4491.864
$ws.Range("I31").Value = 2858.2
$ws.Range("J31").Value = 5853.25
$ws.Range("K31").Value = 2858.2
$ws.Range("L31").Value = 5853.25
$ws.Range("M31").Value = -2563.2
$ws.Range("N31").Value = -6443.25
$ws.Range("H34").Value = 4491.864
$ws.Range("I34").Value = 2858.2
$ws.Range("J34").Value = 5853.25
$ws.Range("K34").Value = 2858.2
$ws.Range("L34").Value = 5853.25
$ws.Range("M34").Value = -2656.2
$ws.Range("N34").Value = -6257.25
$ws.Range("H62").Value = 69047.164
$ws.Range("I62").Value = 2856.8
$ws.Range("K62").Value = 2856.8
$ws.Range("M62").Value = -2232.8
$ws.Range("H65").Value = 69047.164
$ws.Range("I65").Value = 2856.8
$ws.Range("K65").Value = 14284
$ws.Range("M65").Value = -11164
$ws.Range("H86").Value = 8891.714
$ws.Range("I86").Value = 4822.875
$ws.Range("K86").Value = 4822.875
$ws.Range("M86").Value = -3699.875
$ws.Range("H89").Value = 8891.714
$ws.Range("I89").Value = 4822.875
$ws.Range("K89").Value = 24114.375
$ws.Range("M89").Value = -18498.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 7692390
$ws.Range("I7").Value = 10000057
$ws.Range("J7").Value = 166.33333
$ws.Range("K7").Value = 30000171
$ws.Range("L7").Value = 498.99999
$ws.Range("M7").Value = -30000059
$ws.Range("N7").Value = -722.99999
$ws.Range("H70").Value = 2444
$ws.Range("H73").Value = 2444
$ws.Range("H120").Value = 11297.667
$ws.Range("I120").Value = 4515
$ws.Range("K120").Value = 13545
$ws.Range("M120").Value = -8707
$ws.Range("H139").Value = 3289.3635
$ws.Range("I139").Value = 3289.3635
$ws.Range("K139").Value = 9868.0905
$ws.Range("M139").Value = -4728.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 319
$ws.Range("I2").Value = 89
$ws.Range("J2").Value = 1009
$ws.Range("K2").Value = 89
$ws.Range("L2").Value = 1009
$ws.Range("M2").Value = 24
$ws.Range("N2").Value = -1235
$ws.Range("H12").Value = 8900000
$ws.Range("J12").Value = 1500000
$ws.Range("L12").Value = 1500000
$ws.Range("N12").Value = -1500280
$ws.Range("H141").Value = 62049.668
$ws.Range("J141").Value = 62049.668
$ws.Range("L141").Value = 62049.668
$ws.Range("N141").Value = -72409.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents() | Out-Null
$ws.Range("H61").Value = 2794.0715
$ws.Range("I61").Value = 2374.0908
$ws.Range("K61").Value = 2374.0908
$ws.Range("M61").Value = -2172.0908
$ws.Range("H68").Value = 2999.5
$ws.Range("J68").Value = 2999.5
$ws.Range("L68").Value = 2999.5
$ws.Range("N68").Value = -4497.5
$ws.Range("H71").Value = 2999.5
$ws.Range("J71").Value = 2999.5
$ws.Range("L71").Value = 14997.5
$ws.Range("N71").Value = -22485.5
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340
$ws.Range("H113").Value = 2794.0715
$ws.Range("I113").Value = 2374.0908
$ws.Range("K113").Value = 2374.0908
$ws.Range("M113").Value = -204.0907999999999
$ws.Range("H132").Value = 5263.1143
$ws.Range("I132").Value = 4898.909
$ws.Range("K132").Value = 14696.727
$ws.Range("M132").Value = -12166.727
$ws.Range("H136").Value = 3533.353
$ws.Range("I136").Value = 3397.7856
$ws.Range("J136").Value = 4166
$ws.Range("K136").Value = 10193.3568
$ws.Range("L136").Value = 12498
$ws.Range("M136").Value = -7643.356800000001
$ws.Range("N136").Value = -17598

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 25000
$ws.Range("J34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("N34").Value = -25406
$ws.Range("H39").Value = 30000
$ws.Range("J39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30826
$ws.Range("H42").Value = 20000
$ws.Range("J42").Value = 20000
$ws.Range("L42").Value = 20000
$ws.Range("N42").Value = -20756
$ws.Range("H62").Value = 6501.1113
$ws.Range("I62").Value = 1007.5
$ws.Range("J62").Value = 8070.7144
$ws.Range("K62").Value = 1007.5
$ws.Range("L62").Value = 8070.7144
$ws.Range("M62").Value = -383.5
$ws.Range("N62").Value = -9318.714400000001
$ws.Range("H65").Value = 6501.1113
$ws.Range("I65").Value = 1007.5
$ws.Range("J65").Value = 8070.7144
$ws.Range("K65").Value = 5037.5
$ws.Range("L65").Value = 40353.572
$ws.Range("M65").Value = -1917.5
$ws.Range("N65").Value = -46593.572
$ws.Range("H125").Value = 66899.5
$ws.Range("J125").Value = 66899.5
$ws.Range("L125").Value = 66899.5
$ws.Range("N125").Value = -76739.5
$ws.Range("H132").Value = 1611.7
$ws.Range("I132").Value = 1608.5
$ws.Range("K132").Value = 4825.5
$ws.Range("M132").Value = -2295.5
$ws.Range("H136").Value = 57284.555
$ws.Range("I136").Value = 1090
$ws.Range("J136").Value = 203390.4
$ws.Range("K136").Value = 3270
$ws.Range("L136").Value = 610171.2
$ws.Range("M136").Value = -720
$ws.Range("N136").Value = -615271.2
$ws.Range("H140").Value = 97450
$ws.Range("J140").Value = 97450
$ws.Range("L140").Value = 97450
$ws.Range("N140").Value = -107810
